$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark rows 2-9 in column B ("Fixat?") as "Ja" (were "Nej")
$ws.Range("B2:B9").Value = "Ja"

# Move active selection to A18 (row for "Meny")
$ws.Range("A18").Select()
